$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D17 from "Price" to "Quantity" (same value currently found in D18)
$ws.Range("D17").Value = $ws.Range("D18").Text

# Delete row 18 entirely; rows below (19-25) shift up by one
$ws.Rows("18").Delete()

# Update the active cell selection to match the saved view state
$ws.Range("G23").Select()
